# Update the "取得日時" (retrieved-at) timestamp column on the active
# sheet ("ランサーズ") for every existing data row (A2:A21) to the new
# scrape time: 2025-10-29 01:53:35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-29 01:53:35"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
